$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:K1
$ws.Range("F1").Value = "C4.5 acc"
$ws.Range("G1").Value = "credal-C4.5 acc"
$ws.Range("H1").Value = "SPN acc"
$ws.Range("I1").Value = "CSPN low"
$ws.Range("J1").Value = "CSPN high"
$ws.Range("K1").Value = "CSPN robust"

# Copy header style (bold, border, centered) from A1 to new headers F1:K1
$ws.Range("A1").Copy()
$ws.Range("F1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update existing B:E values and add new F:K values for rows 2-6
# Row 2
$ws.Range("B2").Value = 72.09876543209877
$ws.Range("C2").Value = 63.95061728395061
$ws.Range("D2").Value = 99.01234567901234
$ws.Range("E2").Value = 98.48043454101415
$ws.Range("F2").Value = 90.12345679012347
$ws.Range("G2").Value = 88.76543209876543
$ws.Range("H2").Value = 88.08641975308642
$ws.Range("I2").Value = 88.08641975308642
$ws.Range("J2").Value = 88.08641975308642
$ws.Range("K2").Value = 88.08641975308642

# Row 3
$ws.Range("B3").Value = 72.71604938271604
$ws.Range("C3").Value = 63.14814814814815
$ws.Range("D3").Value = 97.77777777777779
$ws.Range("E3").Value = 96.58919010136802
$ws.Range("F3").Value = 86.35802469135803
$ws.Range("G3").Value = 87.71604938271604
$ws.Range("H3").Value = 88.39506172839506
$ws.Range("I3").Value = 88.33333333333333
$ws.Range("J3").Value = 88.33333333333333
$ws.Range("K3").Value = 88.33333333333333

# Row 4
$ws.Range("B4").Value = 69.19753086419753
$ws.Range("C4").Value = 60.8641975308642
$ws.Range("D4").Value = 98.88888888888889
$ws.Range("E4").Value = 98.17355444463996
$ws.Range("F4").Value = 84.19753086419753
$ws.Range("G4").Value = 87.34567901234568
$ws.Range("H4").Value = 87.8395061728395
$ws.Range("I4").Value = 88.14814814814814
$ws.Range("J4").Value = 88.14814814814814
$ws.Range("K4").Value = 88.14814814814814

# Row 5
$ws.Range("B5").Value = 68.27160493827161
$ws.Range("C5").Value = 59.5679012345679
$ws.Range("D5").Value = 98.51851851851852
$ws.Range("E5").Value = 97.45465402528806
$ws.Range("F5").Value = 82.90123456790123
$ws.Range("G5").Value = 86.35802469135803
$ws.Range("H5").Value = 86.54320987654322
$ws.Range("I5").Value = 86.35802469135803
$ws.Range("J5").Value = 86.48148148148148
$ws.Range("K5").Value = 86.4651758676916

# Row 6
$ws.Range("B6").Value = 69.50617283950618
$ws.Range("C6").Value = 57.16049382716049
$ws.Range("D6").Value = 98.76543209876543
$ws.Range("E6").Value = 97.85453499433001
$ws.Range("F6").Value = 81.85185185185186
$ws.Range("G6").Value = 85.37037037037038
$ws.Range("H6").Value = 87.5925925925926
$ws.Range("I6").Value = 87.22222222222223
$ws.Range("J6").Value = 87.53086419753087
$ws.Range("K6").Value = 87.48427672955975
